# Apply "Added few more changes" update to Sheet1 of the Paypal Payment
# Details workbook:
#   - D17:D20 label text "FIXED" -> "FIXED - DONE" (keep existing style/fill,
#     but the text is now shown in green bold)
#   - D21 label text "DONE" -> "DONE - CC" (green bold, no fill)
#   - New row 22: A22 = "Home Needs", B22 = 5000, C22 = running total formula
#   - Active cell selection moves from F22 to D22

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update STATUS column labels (rows 17-20 share one wording, row 21 differs)
$ws.Range("D17").Value = "FIXED - DONE"
$ws.Range("D18").Value = "FIXED - DONE"
$ws.Range("D19").Value = "FIXED - DONE"
$ws.Range("D20").Value = "FIXED - DONE"

# --- Fill in the new "Home Needs" row before renaming D21, so the freed up
#     shared-string slot ("DONE") is reused for "Home Needs" first, matching
#     the original authoring order.
$ws.Range("A22").Value = "Home Needs"
$ws.Range("B22").Value = 5000
$ws.Range("C22").Formula = "=C21-B22"

$ws.Range("D21").Value = "DONE - CC"

# --- Recolor the STATUS labels green (bold) to reflect the fixed/done state
$ws.Range("D17:D21").Font.Color = 5287936

# --- Move the active selection to D22, matching the saved view state
$ws.Activate()
$ws.Range("D22").Select()
